$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "1725236280_w_verifier"
$ws.Range("A5").Value = "1725236280_wo_verifier"
$ws.Range("A2").Value = "1724982582_w_verifier"
$ws.Range("A3").Value = "1724982582_wo_verifier"
$ws.Range("B2").Value = "pass"
$ws.Range("B3").Value = "fail"
$ws.Range("B4").Value = "pass"
$ws.Range("B5").Value = "fail"

$ws.Columns.Item(1).ColumnWidth = 21.54296875
$ws.Columns.Item(2).ColumnWidth = 5.90625

$ws.Range("J8").Select()
